# Append: 2026-01-22 18:39 JST
# Update the "取得日時" (acquired timestamp) column on the "ランサーズ" sheet
# for the existing data rows (2-8) from 2026-01-22 18:29:40 to 2026-01-22 18:39:08.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = "2026-01-22 18:39:08"
}
